$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 501, shifting existing rows 501:569 down to 502:570
$ws.Rows.Item(501).Insert()

# Copy the number format (date style) from the row above (row 500, column D)
# into the new row's D cell so it keeps the same date formatting.
$ws.Range("D500").Copy()
$ws.Range("D501").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 501 with the latest weekly data point
$ws.Cells.Item(501, 1).Value = 8
$ws.Cells.Item(501, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(501, 3).Value = "Coquimbo"
$ws.Cells.Item(501, 4).Value = 45127
$ws.Cells.Item(501, 5).Value = 4
$ws.Cells.Item(501, 6).Value = 100114013
$ws.Cells.Item(501, 7).Value = "Zanahoria"
$ws.Cells.Item(501, 8).Value = "Sin especificar"
$ws.Cells.Item(501, 9).Value = "Primera"
$ws.Cells.Item(501, 10).Value = 400
$ws.Cells.Item(501, 11).Value = 5000
$ws.Cells.Item(501, 12).Value = 6000
$ws.Cells.Item(501, 13).Value = 5500
$ws.Cells.Item(501, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(501, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(501, 16).Value = 275
$ws.Cells.Item(501, 17).Value = 20
$ws.Cells.Item(501, 18).Value = "Hortaliza"
